$wb = $excel.ActiveWorkbook

# --- Rename worksheets ---
# Before:  1=inflow  2=inflow2  3=outflow  4=bypass
# After:   1=inflow1 2=inflow2  3=bypass   4=outflow
# (sheetId/r:id stay paired with their position; only display names
#  at positions 1, 3 and 4 change - "outflow" and "bypass" swap names)

# Use a temporary name to avoid a name collision while swapping 3 & 4.
$wb.Worksheets.Item(3).Name = "outflow_tmp"
$wb.Worksheets.Item(4).Name = "outflow"
$wb.Worksheets.Item("outflow_tmp").Name = "bypass"

$wb.Worksheets.Item(1).Name = "inflow1"

# --- Update selections / active sheet to match the recorded session ---
# The sheet now called "outflow" (position 4) ends up with F28 selected
# and is no longer the active tab.
$wsOut = $wb.Worksheets.Item("outflow")
$wsOut.Activate()
$wsOut.Range("F28").Select()

# The sheet now called "inflow1" (position 1) ends up with C30 selected
# and becomes the active tab (selected last).
$wsIn1 = $wb.Worksheets.Item("inflow1")
$wsIn1.Activate()
$wsIn1.Range("C30").Select()
